# Applies the "changed the name and price" commit to the workbook.
# Sheet1 (Table1 data rows) gets new category/SKU/name/image/seller text plus
# refreshed price numbers; Sheet2 (variant price list) gets refreshed
# name/price text and numbers too. Finally the selections on both sheets are
# updated to match the author's last-saved cursor position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 ("Sheet1") - the Table1 product rows (A2:AA4)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")

# Row 2 (SKU01)
$ws1.Range("A2").Value = "Category1"
$ws1.Range("B2").Value = "SKU01"
$ws1.Range("C2").Value = "ImportTest"
$ws1.Range("D2").Value = "ImportTest"
$ws1.Range("E2").Value = "ImportTest"
$ws1.Range("F2").Value = "ImportTest.jpg"
$ws1.Range("G2").Value = "ImportTest_thumb.jpg"
$ws1.Range("J2").Value = "New"
$ws1.Range("K2").Value = "Seller detail location and contact and google map description here"
$ws1.Range("L2").Value = "USD"
$ws1.Range("M2").Value = 100
$ws1.Range("N2").Value = 80
$ws1.Range("O2").Value = 125
$ws1.Range("Q2").Value = "ecommerce, online shopping, nepal"
$ws1.Range("R2").Value = "ecommerce, online shopping, nepal"
$ws1.Range("S2").Value = "ecommerce, online shopping, nepal"

# Row 3 (SKU02)
$ws1.Range("A3").Value = "Category1/SubCategory/ChildCategory"
$ws1.Range("B3").Value = "SKU02"
$ws1.Range("C3").Value = "ImportTest1"
$ws1.Range("D3").Value = "ImportTest1"
$ws1.Range("E3").Value = "ImportTest1"
$ws1.Range("F3").Value = "ImportTest1.png"
$ws1.Range("G3").Value = "ImportTest1_thumb.png"
$ws1.Range("H3").Value = 10
$ws1.Range("J3").Value = "Used"
$ws1.Range("K3").Value = "Seller detail location and contact and google map description here"
$ws1.Range("L3").Value = "USD"
$ws1.Range("M3").Value = 150
$ws1.Range("N3").Value = 100
$ws1.Range("O3").Value = 180
$ws1.Range("P3").Value = 120
$ws1.Range("Q3").Value = "ecommerce, online shopping, nepal"
$ws1.Range("R3").Value = "ecommerce, online shopping, nepal"
$ws1.Range("S3").Value = "ecommerce, online shopping, nepal"

# Row 4 (SKU03)
$ws1.Range("A4").Value = "Category2"
$ws1.Range("B4").Value = "SKU03"
$ws1.Range("C4").Value = "ImportTest2"
$ws1.Range("D4").Value = "ImportTest2"
$ws1.Range("E4").Value = "ImportTest2"
$ws1.Range("F4").Value = "ImportTest2.jpg"
$ws1.Range("G4").Value = "ImportTest2_thumb.jpg"
$ws1.Range("H4").Value = 5
$ws1.Range("J4").Value = "Not Specified"
$ws1.Range("K4").Value = "Seller detail location and contact and google map description here"
$ws1.Range("L4").Value = "NPR"
$ws1.Range("Q4").Value = "ecommerce, online shopping, nepalsss"
$ws1.Range("R4").Value = "ecommerce, online shopping, nepalsss"
$ws1.Range("S4").Value = "ecommerce, online shopping, nepalsss"
$ws1.Range("U4").Value = 5
$ws1.Range("Y4").Value = 95

# ---------------------------------------------------------------------------
# Sheet2 ("Sheet2") - per-variant price list (A2:H5)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")

# Row 2
$ws2.Range("A2").Value = "SKU01"
$ws2.Range("B2").Value = 5
$ws2.Range("D2").Value = 105
$ws2.Range("E2").Value = "Color"
$ws2.Range("F2").Value = "red"
$ws2.Range("H2").ClearContents()

# Row 3
$ws2.Range("A3").Value = "SKU01"
$ws2.Range("D3").Value = 95
$ws2.Range("E3").Value = "Color"
$ws2.Range("F3").Value = "green"

# Row 4
$ws2.Range("A4").Value = "SKU02"
$ws2.Range("D4").Value = 145
$ws2.Range("E4").Value = "size"
$ws2.Range("F4").Value = "30ml"

# Row 5
$ws2.Range("A5").Value = "SKU02"
$ws2.Range("E5").Value = "size"
$ws2.Range("F5").Value = "60ml"

# ---------------------------------------------------------------------------
# Restore the author's last selection on each sheet
# ---------------------------------------------------------------------------
$ws1.Activate()
$ws1.Range("S4").Select()

$ws2.Activate()
$ws2.Range("D5").Select()

$ws1.Activate()
